# rerun the analysis scripts
# Adds two new columns (ci.lower, ci.upper) to the N.z_ml_results sheet,
# populated with the new confidence-interval bounds for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the confidence-interval columns
$ws.Cells.Item(1, 7).Value = "ci.lower"
$ws.Cells.Item(1, 8).Value = "ci.upper"

# Row 2
$ws.Cells.Item(2, 7).Value = -0.423496272408135
$ws.Cells.Item(2, 8).Value = 0.145413167047354

# Row 3
$ws.Cells.Item(3, 7).Value = -0.0955011327947605
$ws.Cells.Item(3, 8).Value = -0.0143313222334534

# Row 4
$ws.Cells.Item(4, 7).Value = -0.0859914221004081
$ws.Cells.Item(4, 8).Value = 0.00488064388526915

# Row 5
$ws.Cells.Item(5, 7).Value = -0.594112850517934
$ws.Cells.Item(5, 8).Value = -0.0891552011441182

# Row 6
$ws.Cells.Item(6, 7).Value = -0.537749767255214
$ws.Cells.Item(6, 8).Value = 0.0305212432734801

# Row 7
$ws.Cells.Item(7, 7).Value = -0.59566188066955
$ws.Cells.Item(7, 8).Value = -0.0893876554575121

# Row 8
$ws.Cells.Item(8, 7).Value = -0.536347692543689
$ws.Cells.Item(8, 8).Value = 0.0304416652504594

# Row 9
$ws.Cells.Item(9, 7).Value = -0.0882270134932541
$ws.Cells.Item(9, 8).Value = -0.00724460312842229

# Row 10
$ws.Cells.Item(10, 7).Value = -0.353328388874234
$ws.Cells.Item(10, 8).Value = -0.294465346905277

# Row 11
$ws.Cells.Item(11, 7).Value = -0.0437406042767956
$ws.Cells.Item(11, 8).Value = 0.0150189274637207

# Row 12
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(13, 8).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(14, 8).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 7).Value = -0.0874083900338857
$ws.Cells.Item(15, 8).Value = 0.0206584502252841

# Row 16
$ws.Cells.Item(16, 7).Value = -0.545185637805737
$ws.Cells.Item(16, 8).Value = 0.128851364929424

# Row 17
$ws.Cells.Item(17, 7).Value = -0.547899870805125
$ws.Cells.Item(17, 8).Value = 0.128691110663551

# Row 18
$ws.Cells.Item(18, 7).Value = -0.171982844200816
$ws.Cells.Item(18, 8).Value = 0.00976128777053829

# Row 19
$ws.Cells.Item(19, 7).Value = -1.07269538508738
$ws.Cells.Item(19, 8).Value = 0.0608833305009188

# Row 20
$ws.Cells.Item(20, 7).Value = -1.07549953451043
$ws.Cells.Item(20, 8).Value = 0.0610424865469602

# Row 21
$ws.Cells.Item(21, 7).Value = -0.0150189274637207
$ws.Cells.Item(21, 8).Value = 0.0437406042767956

# Row 22
$ws.Cells.Item(22, 7).Value = 0.0144892062568446
$ws.Cells.Item(22, 8).Value = 0.176454026986508

# Row 23
$ws.Cells.Item(23, 7).Value = -0.0936764027491226
$ws.Cells.Item(23, 8).Value = 0.272819911582955

# Row 24
$ws.Cells.Item(24, 7).Value = 0.0903724133504158
$ws.Cells.Item(24, 8).Value = 1.10058315006988

# Row 25
$ws.Cells.Item(25, 7).Value = -0.0954611823506244
$ws.Cells.Item(25, 8).Value = 0.271500710030943

# Row 26
$ws.Cells.Item(26, 7).Value = 0.0900552521092615
$ws.Cells.Item(26, 8).Value = 1.10044132353452
